$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row at the top; everything currently on rows 1..23 shifts
# down to rows 2..24 (keeping its existing content/types untouched).
$ws.Rows.Item(1).Insert()

# Row 2 now holds what used to be row 1 ("bernat"). Update it in place:
# the name was corrected to "xernat" and the associated counters changed.
$ws.Cells.Item(2, 2).Value = "xernat"
$ws.Cells.Item(2, 6).Value = 0
$ws.Cells.Item(2, 7).Value = 26
$ws.Cells.Item(2, 8).Value = 8

# Column L on row 2 flips from "true" to "false". Typing the word "false"
# directly would be auto-converted to an Excel boolean, so instead copy the
# literal text "false" that already lives in column I of the same row.
$ws.Cells.Item(2, 9).Copy($ws.Cells.Item(2, 12))
